$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    "A1" = 2204.2292170912619;  "B1" = 1383.0965219650079;  "C1" = 1428.823605990543
    "A2" = 2227.1192368922416;  "B2" = 1482.0442398488944;  "C2" = 1329.6420023845887
    "A3" = 2343.9137069531039;  "B3" = 1585.8750567376301;  "C3" = 1454.0892310244092
    "A4" = 2320.9625073864663;  "B4" = 1788.7195930546336;  "C4" = 1763.2947254090432
    "A5" = 2422.6539933955587;  "B5" = 1669.0272781903411;  "C5" = 1624.4305377220869
    "A6" = 2360.7084402309279;  "B6" = 1774.6285770614304;  "C6" = 1784.432976973259
    "A7" = 1992.9848472507151;  "B7" = 1566.2947814339752;  "C7" = 1482.2945271321751
    "A8" = 2135.6723079001295;  "B8" = 1612.0574530051363;  "C8" = 1679.2576002754586
    "A9" = 2471.0153725349251;  "B9" = 1788.3095594075089;  "C9" = 1513.8952422040152
    "A10" = 2111.5947501270725; "B10" = 1359.5436219266949; "C10" = 1316.8044372235104
    "A11" = 1970.003883917462;  "B11" = 1416.252990930532;  "C11" = 1298.4013269022387
    "A12" = 2787.953718940531;  "B12" = 2270.4393957721932; "C12" = 2036.7347968306058
    "A13" = 2315.1853419622416; "B13" = 1782.1432978595883; "C13" = 1792.043876782089
    "A14" = 2593.0484713468654; "B14" = 1922.7160414204643; "C14" = 1702.3625413241564
    "A15" = 2508.3188497746851; "B15" = 2002.6740676254963; "C15" = 2042.0927244971056
    "A16" = 2205.5211166450963; "B16" = 1510.1431157823199; "C16" = 1270.5712433044525
    "A17" = 2225.5452594753074; "B17" = 1686.5873845206668; "C17" = 1564.3984569351537
    "A18" = 2487.9072684430239; "B18" = 2061.2592440794292; "C18" = 1916.9784886481884
    "A19" = 1742.2116480415741; "B19" = 1927.5181621920403; "C19" = 1884.7058218307263
    "A20" = 2351.256836519733;  "B20" = 1850.0736740607986; "C20" = 1653.3462287517611
    "A21" = 2583.1100843896215; "B21" = 1907.0148737384604; "C21" = 1813.0877816457973
    "A22" = 2448.3691808016056; "B22" = 1890.2457143793833; "C22" = 1644.281117503079
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
